$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The used range (A1:J7) is made up of four dense 4x2 blocks (no internal
# gaps), each of which needs to move down 1 row and right 2 columns:
#   A1:D2 -> C2:F3
#   G1:J2 -> I2:L3
#   A6:D7 -> C7:F8
#   G6:J7 -> I7:L8
# Copy (not Cut) each block to its destination so the original values stay
# intact until every block has been duplicated; Cut/Copy of the full sparse
# A1:J7 range in one shot does not handle the internal gaps correctly, so
# the range is split into its four contiguous rectangles instead.
$ws.Range("A1:D2").Copy($ws.Range("C2"))
$ws.Range("G1:J2").Copy($ws.Range("I2"))
$ws.Range("A6:D7").Copy($ws.Range("C7"))
$ws.Range("G6:J7").Copy($ws.Range("I7"))

# Now remove the leftover original cells that are not part of the new
# C2:L8 block (i.e. the parts of A1:J7 that don't overlap the destination).
$ws.Range("A1:J1").ClearContents()
$ws.Range("A2:B2").ClearContents()
$ws.Range("G2:H2").ClearContents()
$ws.Range("A6:J6").ClearContents()
$ws.Range("A7:B7").ClearContents()
$ws.Range("G7:H7").ClearContents()
